# Contest 8 PBKS vs CSK
# Fill in the predictor point inputs for row 17 (Contest 8, "PBKS vs CSK"),
# which lets the existing VLOOKUP/RANK formulas in D17,G17,J17,M17,P17,S17
# (and the SUM totals in row 27) calculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E17").Value = 80
$ws.Range("H17").Value = 20
$ws.Range("K17").Value = 100
$ws.Range("N17").Value = 0
$ws.Range("Q17").Value = 40
$ws.Range("T17").Value = 60

$excel.Calculate()
